# Npnt -> Itgb1 LR-pair sheet: refresh with updated TPM-derived expression data.
#
# The underlying per-cluster ligand (Npnt) and receptor (Itgb1) TPM inputs were
# recomputed; every downstream NATMI column (detection rate, average/total
# expression, derived specificity, and edge weights/specificities) is refreshed
# below with the newly computed values for each Sending/Target cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster = ECs, Target cluster = ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.660320666666667
$ws.Range("H2").Value = 4.980962
$ws.Range("I2").Value = 0.3342448133445559
$ws.Range("J2").Value = 0.3342448133445559
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 241.9099575187898
$ws.Range("R2").Value = 2177.189617669108
$ws.Range("S2").Value = 0.09579248283412437
$ws.Range("T2").Value = 0.09579248283412437

# Row 3: Sending cluster = ECs, Target cluster = FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.660320666666667
$ws.Range("H3").Value = 4.980962
$ws.Range("I3").Value = 0.3342448133445559
$ws.Range("J3").Value = 0.3342448133445559
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 280.2616348646552
$ws.Range("R3").Value = 2522.354713781896
$ws.Range("S3").Value = 0.1109791350558641
$ws.Range("T3").Value = 0.1109791350558641

# Row 4: Sending cluster = ECs, Target cluster = MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.660320666666667
$ws.Range("H4").Value = 4.980962
$ws.Range("I4").Value = 0.3342448133445559
$ws.Range("J4").Value = 0.3342448133445559
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 212.7305025335965
$ws.Range("R4").Value = 1914.574522802368
$ws.Range("S4").Value = 0.08423788429900153
$ws.Range("T4").Value = 0.08423788429900153

# Row 5: Sending cluster = ECs, Target cluster = Resolving-Mac
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.660320666666667
$ws.Range("H5").Value = 4.980962
$ws.Range("I5").Value = 0.3342448133445559
$ws.Range("J5").Value = 0.3342448133445559
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 109.1844785259993
$ws.Range("R5").Value = 982.6603067339939
$ws.Range("S5").Value = 0.04323531115556593
$ws.Range("T5").Value = 0.04323531115556593

# Row 6: Sending cluster = FAPs, Target cluster = ECs
$ws.Range("I6").Value = 0.02135559483851209
$ws.Range("J6").Value = 0.02135559483851209
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 15.45612926189956
$ws.Range("R6").Value = 139.105163357096
$ws.Range("S6").Value = 0.006120380542365728
$ws.Range("T6").Value = 0.006120380542365728

# Row 7: Sending cluster = FAPs, Target cluster = FAPs
$ws.Range("I7").Value = 0.02135559483851209
$ws.Range("J7").Value = 0.02135559483851209
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("Q7").Value = 17.90649752515023
$ws.Range("S7").Value = 0.007090687272201315
$ws.Range("T7").Value = 0.007090687272201315

# Row 8: Sending cluster = FAPs, Target cluster = MuSCs
$ws.Range("I8").Value = 0.02135559483851209
$ws.Range("J8").Value = 0.02135559483851209
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 13.59179332191289
$ws.Range("R8").Value = 122.326139897216
$ws.Range("S8").Value = 0.005382133260774013
$ws.Range("T8").Value = 0.005382133260774013

# Row 9: Sending cluster = FAPs, Target cluster = Resolving-Mac
$ws.Range("I9").Value = 0.02135559483851209
$ws.Range("J9").Value = 0.02135559483851209
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 6.976022941758667
$ws.Range("R9").Value = 62.784206475828
$ws.Range("S9").Value = 0.002762393763171035
$ws.Range("T9").Value = 0.002762393763171035

# Row 10: Sending cluster = MuSCs, Target cluster = ECs
$ws.Range("G10").Value = 3.108009
$ws.Range("H10").Value = 9.324027000000001
$ws.Range("I10").Value = 0.6256838868143543
$ws.Range("J10").Value = 0.6256838868143542
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 452.8392257307021
$ws.Range("R10").Value = 4075.553031576319
$ws.Range("S10").Value = 0.1793171070854209
$ws.Range("T10").Value = 0.1793171070854209

# Row 11: Sending cluster = MuSCs, Target cluster = FAPs
$ws.Range("G11").Value = 3.108009
$ws.Range("H11").Value = 9.324027000000001
$ws.Range("I11").Value = 0.6256838868143543
$ws.Range("J11").Value = 0.6256838868143542
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 524.6309950853241
$ws.Range("R11").Value = 4721.678955767917
$ws.Range("S11").Value = 0.2077455021133514
$ws.Range("T11").Value = 0.2077455021133514

# Row 12: Sending cluster = MuSCs, Target cluster = MuSCs
$ws.Range("G12").Value = 3.108009
$ws.Range("H12").Value = 9.324027000000001
$ws.Range("I12").Value = 0.6256838868143543
$ws.Range("J12").Value = 0.6256838868143542
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 398.2172418393921
$ws.Range("R12").Value = 3583.955176554528
$ws.Range("S12").Value = 0.157687673109485
$ws.Range("T12").Value = 0.1576876731094849

# Row 13: Sending cluster = MuSCs, Target cluster = Resolving-Mac
$ws.Range("G13").Value = 3.108009
$ws.Range("H13").Value = 9.324027000000001
$ws.Range("I13").Value = 0.6256838868143543
$ws.Range("J13").Value = 0.6256838868143542
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 204.386025381711
$ws.Range("R13").Value = 1839.474228435399
$ws.Range("S13").Value = 0.080933604506097
$ws.Range("T13").Value = 0.08093360450609699

# Row 14: Sending cluster = Resolving-Mac, Target cluster = ECs
$ws.Range("G14").Value = 0.09296800000000001
$ws.Range("H14").Value = 0.278904
$ws.Range("I14").Value = 0.01871570500257782
$ws.Range("J14").Value = 0.01871570500257782
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 13.54550683017067
$ws.Range("R14").Value = 121.909561471536
$ws.Range("S14").Value = 0.005363804548673253
$ws.Range("T14").Value = 0.005363804548673253

# Row 15: Sending cluster = Resolving-Mac, Target cluster = FAPs
$ws.Range("G15").Value = 0.09296800000000001
$ws.Range("H15").Value = 0.278904
$ws.Range("I15").Value = 0.01871570500257782
$ws.Range("J15").Value = 0.01871570500257782
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 15.69297075751467
$ws.Range("R15").Value = 141.236736817632
$ws.Range("S15").Value = 0.006214165995167342
$ws.Range("T15").Value = 0.006214165995167342

# Row 16: Sending cluster = Resolving-Mac, Target cluster = MuSCs
$ws.Range("G16").Value = 0.09296800000000001
$ws.Range("H16").Value = 0.278904
$ws.Range("I16").Value = 0.01871570500257782
$ws.Range("J16").Value = 0.01871570500257782
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 11.91163234705067
$ws.Range("R16").Value = 107.204691123456
$ws.Range("S16").Value = 0.004716816326349955
$ws.Range("T16").Value = 0.004716816326349955

# Row 17: Sending cluster = Resolving-Mac, Target cluster = Resolving-Mac
$ws.Range("G17").Value = 0.09296800000000001
$ws.Range("H17").Value = 0.278904
$ws.Range("I17").Value = 0.01871570500257782
$ws.Range("J17").Value = 0.01871570500257782
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 6.113675992472
$ws.Range("R17").Value = 55.02308393224801
$ws.Range("S17").Value = 0.00242091813238727
$ws.Range("T17").Value = 0.00242091813238727
